# Update UnitPriceUSD (I) and UnitPriceEURO (J) columns with refreshed
# currency-conversion values for rows 2-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; I = 44.53;   J = 46.81 },
    @{ Row = 3; I = 1099.97; J = 1156.18 },
    @{ Row = 4; I = 298.37;  J = 313.62 },
    @{ Row = 5; I = 133.6;   J = 140.43 },
    @{ Row = 6; I = 2395.88; J = 2518.32 },
    @{ Row = 7; I = 623.46;  J = 655.33 },
    @{ Row = 8; I = 1202.39; J = 1263.84 }
)

foreach ($u in $updates) {
    $ws.Range("I$($u.Row)").Value = $u.I
    $ws.Range("J$($u.Row)").Value = $u.J
}
